# Finished Week 13 logging
# Update the "H" (home) row of target-depth counting stats on both the
# OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet (Short Att, Short Comp, Deep Att, Deep Comp, Short Int, Deep Int) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 135
$wsOff.Range("C2").Value = 99
$wsOff.Range("D2").Value = 39
$wsOff.Range("E2").Value = 20
$wsOff.Range("F2").Value = 2
$wsOff.Range("G2").Value = 2

# --- DEF sheet (Short Att, Short Comp, Deep Att, Deep Comp, Short Int, Deep Int) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 271
$wsDef.Range("C2").Value = 190
$wsDef.Range("D2").Value = 49
$wsDef.Range("E2").Value = 22
$wsDef.Range("F2").Value = 2
$wsDef.Range("G2").Value = 2
